$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.16212522983551
$ws.Range("B1").Value = 2.120952367782593
$ws.Range("C1").Value = 3.445654153823853
$ws.Range("D1").Value = 3.525397777557373
$ws.Range("E1").Value = 1.179910898208618
